$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 850
$ws.Range("I2").Value = 550
$ws.Range("K2").Value = 550
$ws.Range("M2").Value = -437
$ws.Range("H9").Value = 176
$ws.Range("I9").Value = 157.33333
$ws.Range("K9").Value = 157.33333
$ws.Range("M9").Value = 11.66667000000001
$ws.Range("H112").Value = 2137.4783
$ws.Range("I112").Value = 1830.5
$ws.Range("K112").Value = 5491.5
$ws.Range("M112").Value = -4383.5
$ws.Range("H137").Value = 7091.1035
$ws.Range("I137").Value = 1754.4762
$ws.Range("J137").Value = 21099.75
$ws.Range("K137").Value = 5263.4286
$ws.Range("L137").Value = 63299.25
$ws.Range("M137").Value = -2713.4286
$ws.Range("N137").Value = -68399.25

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").ClearContents()
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").ClearContents()
$ws.Range("N65").ClearContents()
$ws.Range("H102").Value = 5589.231
$ws.Range("I102").Value = 5055
$ws.Range("J102").Value = 12000
$ws.Range("K102").Value = 5055
$ws.Range("L102").Value = 12000
$ws.Range("M102").Value = -3433
$ws.Range("N102").Value = -15244
$ws.Range("H110").Value = 1356.25
$ws.Range("I110").Value = 1292.3077
$ws.Range("J110").Value = 1633.3334
$ws.Range("K110").Value = 1292.3077
$ws.Range("L110").Value = 1633.3334
$ws.Range("M110").Value = 752.6922999999999
$ws.Range("N110").Value = -5723.3334

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 12731.889
$ws.Range("I99").Value = 15855.429
$ws.Range("J99").Value = 1799.5
$ws.Range("K99").Value = 15855.429
$ws.Range("L99").Value = 1799.5
$ws.Range("M99").Value = -14357.429
$ws.Range("N99").Value = -4795.5
$ws.Range("H105").Value = 2666.8628
$ws.Range("I105").Value = 2229.5
$ws.Range("J105").Value = 3716.5334
$ws.Range("K105").Value = 2229.5
$ws.Range("L105").Value = 3716.5334
$ws.Range("M105").Value = -482.5
$ws.Range("N105").Value = -7210.5334

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 5414.75
$ws.Range("I99").Value = 5886.3335
$ws.Range("K99").Value = 5886.3335
$ws.Range("M99").Value = -4388.3335
$ws.Range("H126").Value = 5414.75
$ws.Range("I126").Value = 5886.3335
$ws.Range("K126").Value = 17659.0005
$ws.Range("M126").Value = -15189.0005
$ws.Range("H134").Value = 3666.5334
$ws.Range("I134").Value = 2916.5
$ws.Range("K134").Value = 8749.5
$ws.Range("M134").Value = -6214.5

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4667899
$ws.Range("I4").Value = 7000215
$ws.Range("J4").Value = 3267.9
$ws.Range("K4").Value = 21000645
$ws.Range("L4").Value = 9803.700000000001
$ws.Range("M4").Value = -21000533
$ws.Range("N4").Value = -10027.7
$ws.Range("H5").Value = 2312.8965
$ws.Range("I5").Value = 1656.9286
$ws.Range("J5").Value = 2925.1333
$ws.Range("K5").Value = 4970.7858
$ws.Range("L5").Value = 8775.3999
$ws.Range("M5").Value = -4858.7858
$ws.Range("N5").Value = -8999.3999
$ws.Range("H12").Value = 387.57144
$ws.Range("I12").Value = 853
$ws.Range("K12").Value = 2559
$ws.Range("M12").Value = -2386
$ws.Range("H114").Value = 8287.321
$ws.Range("J114").Value = 9843.305
$ws.Range("L114").Value = 29529.915
$ws.Range("N114").Value = -36037.915
$ws.Range("H135").Value = 2312.8965
$ws.Range("I135").Value = 1656.9286
$ws.Range("J135").Value = 2925.1333
$ws.Range("K135").Value = 14912.3574
$ws.Range("L135").Value = 26326.1997
$ws.Range("M135").Value = -12377.3574
$ws.Range("N135").Value = -31396.1997

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 5263286.5
$ws.Range("I2").Value = 6666716
$ws.Range("J2").Value = 424.75
$ws.Range("K2").Value = 6666716
$ws.Range("L2").Value = 424.75
$ws.Range("M2").Value = -6666603
$ws.Range("N2").Value = -650.75
$ws.Range("H70").Value = 66998.336
$ws.Range("J70").Value = 62398
$ws.Range("L70").Value = 62398
$ws.Range("N70").Value = -62938
$ws.Range("H73").Value = 66998.336
$ws.Range("J73").Value = 62398
$ws.Range("L73").Value = 62398
$ws.Range("N73").Value = -64270
$ws.Range("H80").Value = 1594.5
$ws.Range("I80").Value = 1566.1666
$ws.Range("J80").Value = 1679.5
$ws.Range("K80").Value = 1566.1666
$ws.Range("L80").Value = 1679.5
$ws.Range("M80").Value = -568.1666
$ws.Range("N80").Value = -3675.5
$ws.Range("H83").Value = 1594.5
$ws.Range("I83").Value = 1566.1666
$ws.Range("J83").Value = 1679.5
$ws.Range("K83").Value = 7830.833000000001
$ws.Range("L83").Value = 8397.5
$ws.Range("M83").Value = -2838.833000000001
$ws.Range("N83").Value = -18381.5
$ws.Range("H102").Value = 2316.25
$ws.Range("I102").Value = 2316.25
$ws.Range("K102").Value = 2316.25
$ws.Range("M102").Value = -694.25

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5503.4443
$ws.Range("I7").Value = 4218.857
$ws.Range("K7").Value = 4218.857
$ws.Range("M7").Value = -4106.857
$ws.Range("H22").Value = 2070.6775
$ws.Range("J22").Value = 2119.7
$ws.Range("L22").Value = 2119.7
$ws.Range("N22").Value = -2709.7
$ws.Range("H27").Value = 2070.6775
$ws.Range("J27").Value = 2119.7
$ws.Range("L27").Value = 2119.7
$ws.Range("N27").Value = -2333.7
$ws.Range("H40").Value = 3857.7334
$ws.Range("I40").Value = 3572.1667
$ws.Range("K40").Value = 3572.1667
$ws.Range("M40").Value = -3436.1667
$ws.Range("H46").Value = 6557
$ws.Range("J46").Value = 7499.8335
$ws.Range("L46").Value = 7499.8335
$ws.Range("N46").Value = -7875.8335
$ws.Range("H55").Value = 1946.5
$ws.Range("J55").Value = 1727.0526
$ws.Range("L55").Value = 1727.0526
$ws.Range("N55").Value = -2073.0526
$ws.Range("H69").Value = 20000
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 20000
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H74").Value = 41994
$ws.Range("I74").Value = 41994
$ws.Range("K74").Value = 41994
$ws.Range("M74").Value = -40996
$ws.Range("H77").Value = 41994
$ws.Range("I77").Value = 41994
$ws.Range("K77").Value = 125982
$ws.Range("M77").Value = -120990
$ws.Range("H93").Value = 3294.75
$ws.Range("J93").Value = 3903.889
$ws.Range("L93").Value = 3903.889
$ws.Range("N93").Value = -6399.889
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").ClearContents()
$ws.Range("N116").ClearContents()
$ws.Range("H126").Value = 5503.4443
$ws.Range("I126").Value = 4218.857
$ws.Range("K126").Value = 12656.571
$ws.Range("M126").Value = -10186.571
$ws.Range("H136").Value = 26573.8
$ws.Range("I136").Value = 9823
$ws.Range("K136").Value = 29469
$ws.Range("M136").Value = -26919
$ws.Range("H137").Value = 99000
$ws.Range("J137").Value = 99000
$ws.Range("L137").Value = 99000
$ws.Range("N137").Value = -109200
$ws.Range("H139").Value = 42222
$ws.Range("I139").Value = 42222
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 42222
$ws.Range("L139").Value = 0
$ws.Range("M139").ClearContents()
$ws.Range("N139").ClearContents()
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").ClearContents()
$ws.Range("N141").ClearContents()

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 63146.555
$ws.Range("I122").Value = 479.0909
$ws.Range("K122").Value = 1437.2727
$ws.Range("M122").Value = 1012.7273
$ws.Range("H126").Value = 2145.4092
$ws.Range("I126").Value = 2166.0588
$ws.Range("J126").Value = 2075.2
$ws.Range("K126").Value = 6498.176399999999
$ws.Range("L126").Value = 6225.599999999999
$ws.Range("M126").Value = -4028.176399999999
$ws.Range("N126").Value = -11165.6
$ws.Range("H132").Value = 4385.6665
$ws.Range("I132").Value = 3924.4285
$ws.Range("K132").Value = 11773.2855
$ws.Range("M132").Value = -9243.2855
$ws.Range("H136").Value = 4697.8
$ws.Range("I136").Value = 4833.8
$ws.Range("K136").Value = 14501.4
$ws.Range("M136").Value = -11951.4
